# Update the cryptos list table (columns D=Price, E=Volume(1h)) with the
# latest values; also fix row 37/38 ordering swap (RenderToken <-> Stellar).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds free-form text (e.g. "42.966.58", "1.01", "0.0990")
# rather than numeric values, so force text formatting before assigning —
# otherwise Excel auto-converts number-looking strings and can drop
# significant trailing zeros (e.g. "4.50" -> 4.5).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '42.966.58'
$ws.Range("E2").Value = '  -1.85%  '
$ws.Range("D3").Value = '2.246.47'
$ws.Range("E3").Value = '  -1.87%  '
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.36%  '
$ws.Range("D5").Value = '114.64'
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").Value = '299.89'
$ws.Range("E6").Value = '  +12.24%  '
$ws.Range("E7").Value = '  -2.01%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '0.622'
$ws.Range("E9").Value = '  +0.92%  '
$ws.Range("D10").Value = '46.08'
$ws.Range("E10").Value = '  -3.92%  '
$ws.Range("D11").Value = '0.0933'
$ws.Range("E11").Value = '  -0.91%  '
$ws.Range("D12").Value = '56.12'
$ws.Range("E12").Value = '  +2.70%  '
$ws.Range("D13").Value = '9.13'
$ws.Range("E13").Value = '  -0.22%  '
$ws.Range("E14").Value = '  -2.66%  '
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("D16").Value = '0.904'
$ws.Range("E16").Value = '  +1.88%  '
$ws.Range("D17").Value = '2.587.42'
$ws.Range("E17").Value = '  -1.73%  '
$ws.Range("D18").Value = '2.270.66'
$ws.Range("E18").Value = '  -1.02%  '
$ws.Range("D19").Value = '42.926.60'
$ws.Range("E19").Value = '  -1.77%  '
$ws.Range("D20").Value = '7.73'
$ws.Range("E20").Value = '  +11.88%  '
$ws.Range("D22").Value = '3.67'
$ws.Range("E22").Value = '  +26.75%  '
$ws.Range("D23").Value = '73.67'
$ws.Range("E23").Value = '  +1.62%  '
$ws.Range("D24").Value = '2.35'
$ws.Range("E24").Value = '  -4.46%  '
$ws.Range("D25").Value = '232.48'
$ws.Range("E25").Value = '  -1.44%  '
$ws.Range("D26").Value = '9.51'
$ws.Range("E26").Value = '  -0.94%  '
$ws.Range("D27").Value = '12.22'
$ws.Range("E27").Value = '  +4.00%  '
$ws.Range("E28").Value = '  -1.55%  '
$ws.Range("D29").Value = '40.15'
$ws.Range("E29").Value = '  -4.95%  '
$ws.Range("D30").Value = '2.22'
$ws.Range("E30").Value = '  -1.50%  '
$ws.Range("D31").Value = '3.28'
$ws.Range("E31").Value = '  -4.03%  '
$ws.Range("D32").Value = '175.83'
$ws.Range("E32").Value = '  +1.02%  '
$ws.Range("D33").Value = '21.32'
$ws.Range("E33").Value = '  -2.16%  '
$ws.Range("D34").Value = '0.0904'
$ws.Range("E34").Value = '  -1.14%  '
$ws.Range("D35").Value = '5.71'
$ws.Range("E35").Value = '  -0.35%  '
$ws.Range("D36").Value = '4.50'
$ws.Range("E36").Value = '  +14.02%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '4.91'
$ws.Range("E37").Value = '  +4.69%  '
$ws.Range("B38").Value = 'Stellar'
$ws.Range("C38").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D38").Value = '0.129'
$ws.Range("E38").Value = '  -1.40%  '
$ws.Range("D39").Value = '0.0373'
$ws.Range("E39").Value = '  -3.10%  '
$ws.Range("E40").Value = '  -1.46%  '
$ws.Range("E41").Value = '  +1.12%  '
$ws.Range("E42").Value = '  +2.06%  '
$ws.Range("D43").Value = '72.05'
$ws.Range("E43").Value = '  -3.01%  '
$ws.Range("D44").Value = '13.43'
$ws.Range("E44").Value = '  -5.13%  '
$ws.Range("E45").Value = '  +0.34%  '
$ws.Range("D46").Value = '1.34'
$ws.Range("E46").Value = '  -1.84%  '
$ws.Range("E47").Value = '  -5.99%  '
$ws.Range("D48").Value = '1.36'
$ws.Range("E48").Value = '  +6.28%  '
$ws.Range("D49").Value = '106.41'
$ws.Range("E49").Value = '  +4.58%  '
$ws.Range("D50").Value = '8.73'
$ws.Range("E50").Value = '  +1.15%  '
$ws.Range("D51").Value = '0.0990'
$ws.Range("E51").Value = '  -1.41%  '
